# Update L1cam-Ephb2 NATMI edge statistics (columns G:T) with the recomputed TPM values.
# Columns: G=LigAvgExpr H=LigTotExpr I=LigSpecAvg J=LigSpecTot
#          K=RecvCells  L=RecvDetRate M=RecAvgExpr N=RecTotExpr O=RecSpecAvg P=RecSpecTot
#          Q=EdgeAvgWeight R=EdgeTotWeight S=EdgeAvgSpec T=EdgeTotSpec
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 25,14

# Row 2: ECs -> ECs
$data[0,0] = [double]"5.375839"
$data[0,1] = [double]"16.127517"
$data[0,2] = [double]"0.2354568587499626"
$data[0,3] = [double]"0.2354568587499626"
$data[0,4] = 2
$data[0,5] = [double]"0.6666666666666666"
$data[0,6] = [double]"0.2087793333333333"
$data[0,7] = [double]"0.626338"
$data[0,8] = [double]"0.02275344108115409"
$data[0,9] = [double]"0.02275344108115409"
$data[0,10] = [double]"1.122364082527333"
$data[0,11] = [double]"10.101276742746"
$data[0,12] = [double]"0.005357453762720894"
$data[0,13] = [double]"0.005357453762720894"

# Row 3: ECs -> FAPs
$data[1,0] = [double]"5.375839"
$data[1,1] = [double]"16.127517"
$data[1,2] = [double]"0.2354568587499626"
$data[1,3] = [double]"0.2354568587499626"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = [double]"8.445752666666667"
$data[1,7] = [double]"25.337258"
$data[1,8] = [double]"0.9204452022087118"
$data[1,9] = [double]"0.920445202208712"
$data[1,10] = [double]"45.40300656982067"
$data[1,11] = [double]"408.6270591283861"
$data[1,12] = [double]"0.2167251359635374"
$data[1,13] = [double]"0.2167251359635374"

# Row 4: ECs -> Inflammatory-Mac
$data[2,0] = [double]"5.375839"
$data[2,1] = [double]"16.127517"
$data[2,2] = [double]"0.2354568587499626"
$data[2,3] = [double]"0.2354568587499626"
$data[2,4] = 2
$data[2,5] = [double]"0.6666666666666666"
$data[2,6] = [double]"0.05870933333333334"
$data[2,7] = [double]"0.176128"
$data[2,8] = [double]"0.006398331365399365"
$data[2,9] = [double]"0.006398331365399365"
$data[2,10] = [double]"0.3156119237973333"
$data[2,11] = [double]"2.840507314176"
$data[2,12] = [double]"0.001506531004538294"
$data[2,13] = [double]"0.001506531004538294"

# Row 5: ECs -> MuSCs
$data[3,0] = [double]"5.375839"
$data[3,1] = [double]"16.127517"
$data[3,2] = [double]"0.2354568587499626"
$data[3,3] = [double]"0.2354568587499626"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = [double]"0.4594193333333333"
$data[3,7] = [double]"1.378258"
$data[3,8] = [double]"0.05006899181852175"
$data[3,9] = [double]"0.05006899181852175"
$data[3,10] = [double]"2.469764369487333"
$data[3,11] = [double]"22.227879325386"
$data[3,12] = [double]"0.01178908753436671"
$data[3,13] = [double]"0.01178908753436671"

# Row 6: ECs -> Resolving-Mac
$data[4,0] = [double]"5.375839"
$data[4,1] = [double]"16.127517"
$data[4,2] = [double]"0.2354568587499626"
$data[4,3] = [double]"0.2354568587499626"
$data[4,4] = 1
$data[4,5] = [double]"0.3333333333333333"
$data[4,6] = [double]"0.003065"
$data[4,7] = [double]"0.009195"
$data[4,8] = [double]"0.0003340335262130221"
$data[4,9] = [double]"0.0003340335262130222"
$data[4,10] = [double]"0.016476946535"
$data[4,11] = [double]"0.148292518815"
$data[4,12] = [double]"7.865048479929146e-05"
$data[4,13] = [double]"7.865048479929149e-05"

# Row 7: FAPs -> ECs
$data[5,0] = [double]"0.1628146666666667"
$data[5,1] = [double]"0.488444"
$data[5,2] = [double]"0.007131134316291014"
$data[5,3] = [double]"0.007131134316291014"
$data[5,4] = 2
$data[5,5] = [double]"0.6666666666666666"
$data[5,6] = [double]"0.2087793333333333"
$data[5,7] = [double]"0.626338"
$data[5,8] = [double]"0.02275344108115409"
$data[5,9] = [double]"0.02275344108115409"
$data[5,10] = [double]"0.03399233756355555"
$data[5,11] = [double]"0.305931038072"
$data[5,12] = [double]"0.0001622578445075236"
$data[5,13] = [double]"0.0001622578445075236"

# Row 8: FAPs -> FAPs
$data[6,0] = [double]"0.1628146666666667"
$data[6,1] = [double]"0.488444"
$data[6,2] = [double]"0.007131134316291014"
$data[6,3] = [double]"0.007131134316291014"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = [double]"8.445752666666667"
$data[6,7] = [double]"25.337258"
$data[6,8] = [double]"0.9204452022087118"
$data[6,9] = [double]"0.920445202208712"
$data[6,10] = [double]"1.375092405172444"
$data[6,11] = [double]"12.375831646552"
$data[6,12] = [double]"0.006563818367735967"
$data[6,13] = [double]"0.006563818367735967"

# Row 9: FAPs -> Inflammatory-Mac
$data[7,0] = [double]"0.1628146666666667"
$data[7,1] = [double]"0.488444"
$data[7,2] = [double]"0.007131134316291014"
$data[7,3] = [double]"0.007131134316291014"
$data[7,4] = 2
$data[7,5] = [double]"0.6666666666666666"
$data[7,6] = [double]"0.05870933333333334"
$data[7,7] = [double]"0.176128"
$data[7,8] = [double]"0.006398331365399365"
$data[7,9] = [double]"0.006398331365399365"
$data[7,10] = [double]"0.009558740536888889"
$data[7,11] = [double]"0.086028664832"
$data[7,12] = [double]"4.562736036680056e-05"
$data[7,13] = [double]"4.562736036680056e-05"

# Row 10: FAPs -> MuSCs
$data[8,0] = [double]"0.1628146666666667"
$data[8,1] = [double]"0.488444"
$data[8,2] = [double]"0.007131134316291014"
$data[8,3] = [double]"0.007131134316291014"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = [double]"0.4594193333333333"
$data[8,7] = [double]"1.378258"
$data[8,8] = [double]"0.05006899181852175"
$data[8,9] = [double]"0.05006899181852175"
$data[8,10] = [double]"0.0748002056168889"
$data[8,11] = [double]"0.673201850552"
$data[8,12] = [double]"0.0003570487057391545"
$data[8,13] = [double]"0.0003570487057391545"

# Row 11: FAPs -> Resolving-Mac
$data[9,0] = [double]"0.1628146666666667"
$data[9,1] = [double]"0.488444"
$data[9,2] = [double]"0.007131134316291014"
$data[9,3] = [double]"0.007131134316291014"
$data[9,4] = 1
$data[9,5] = [double]"0.3333333333333333"
$data[9,6] = [double]"0.003065"
$data[9,7] = [double]"0.009195"
$data[9,8] = [double]"0.0003340335262130221"
$data[9,9] = [double]"0.0003340335262130222"
$data[9,10] = [double]"0.0004990269533333334"
$data[9,11] = [double]"0.00449124258"
$data[9,12] = [double]"2.382037941569376e-06"
$data[9,13] = [double]"2.382037941569377e-06"

# Row 12: Inflammatory-Mac -> ECs
$data[10,0] = [double]"9.994147"
$data[10,1] = [double]"29.982441"
$data[10,2] = [double]"0.4377345486919088"
$data[10,3] = [double]"0.4377345486919088"
$data[10,4] = 2
$data[10,5] = [double]"0.6666666666666666"
$data[10,6] = [double]"0.2087793333333333"
$data[10,7] = [double]"0.626338"
$data[10,8] = [double]"0.02275344108115409"
$data[10,9] = [double]"0.02275344108115409"
$data[10,10] = [double]"2.086571347895333"
$data[10,11] = [double]"18.779142131058"
$data[10,12] = [double]"0.009959967262846921"
$data[10,13] = [double]"0.009959967262846921"

# Row 13: Inflammatory-Mac -> FAPs
$data[11,0] = [double]"9.994147"
$data[11,1] = [double]"29.982441"
$data[11,2] = [double]"0.4377345486919088"
$data[11,3] = [double]"0.4377345486919088"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = [double]"8.445752666666667"
$data[11,7] = [double]"25.337258"
$data[11,8] = [double]"0.9204452022087118"
$data[11,9] = [double]"0.920445202208712"
$data[11,10] = [double]"84.40809367630867"
$data[11,11] = [double]"759.672843086778"
$data[11,12] = [double]"0.4029106651844632"
$data[11,13] = [double]"0.4029106651844632"

# Row 14: Inflammatory-Mac -> Inflammatory-Mac
$data[12,0] = [double]"9.994147"
$data[12,1] = [double]"29.982441"
$data[12,2] = [double]"0.4377345486919088"
$data[12,3] = [double]"0.4377345486919088"
$data[12,4] = 2
$data[12,5] = [double]"0.6666666666666666"
$data[12,6] = [double]"0.05870933333333334"
$data[12,7] = [double]"0.176128"
$data[12,8] = [double]"0.006398331365399365"
$data[12,9] = [double]"0.006398331365399365"
$data[12,10] = [double]"0.5867497076053334"
$data[12,11] = [double]"5.280747368448"
$data[12,12] = [double]"0.002800770692614376"
$data[12,13] = [double]"0.002800770692614376"

# Row 15: Inflammatory-Mac -> MuSCs
$data[13,0] = [double]"9.994147"
$data[13,1] = [double]"29.982441"
$data[13,2] = [double]"0.4377345486919088"
$data[13,3] = [double]"0.4377345486919088"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = [double]"0.4594193333333333"
$data[13,7] = [double]"1.378258"
$data[13,8] = [double]"0.05006899181852175"
$data[13,9] = [double]"0.05006899181852175"
$data[13,10] = [double]"4.591504351975334"
$data[13,11] = [double]"41.32353916777799"
$data[13,12] = [double]"0.02191692753713949"
$data[13,13] = [double]"0.02191692753713949"

# Row 16: Inflammatory-Mac -> Resolving-Mac
$data[14,0] = [double]"9.994147"
$data[14,1] = [double]"29.982441"
$data[14,2] = [double]"0.4377345486919088"
$data[14,3] = [double]"0.4377345486919088"
$data[14,4] = 1
$data[14,5] = [double]"0.3333333333333333"
$data[14,6] = [double]"0.003065"
$data[14,7] = [double]"0.009195"
$data[14,8] = [double]"0.0003340335262130221"
$data[14,9] = [double]"0.0003340335262130222"
$data[14,10] = [double]"0.030632060555"
$data[14,11] = [double]"0.275688544995"
$data[14,12] = [double]"0.0001462180148448241"
$data[14,13] = [double]"0.0001462180148448241"

# Row 17: MuSCs -> ECs
$data[15,0] = [double]"0.7761303333333333"
$data[15,1] = [double]"2.328391"
$data[15,2] = [double]"0.03399380269149206"
$data[15,3] = [double]"0.03399380269149207"
$data[15,4] = 2
$data[15,5] = [double]"0.6666666666666666"
$data[15,6] = [double]"0.2087793333333333"
$data[15,7] = [double]"0.626338"
$data[15,8] = [double]"0.02275344108115409"
$data[15,9] = [double]"0.02275344108115409"
$data[15,10] = [double]"0.1620399735731111"
$data[15,11] = [double]"1.458359762158"
$data[15,12] = [double]"0.0007734759866652419"
$data[15,13] = [double]"0.000773475986665242"

# Row 18: MuSCs -> FAPs
$data[16,0] = [double]"0.7761303333333333"
$data[16,1] = [double]"2.328391"
$data[16,2] = [double]"0.03399380269149206"
$data[16,3] = [double]"0.03399380269149207"
$data[16,4] = 3
$data[16,5] = 1
$data[16,6] = [double]"8.445752666666667"
$data[16,7] = [double]"25.337258"
$data[16,8] = [double]"0.9204452022087118"
$data[16,9] = [double]"0.920445202208712"
$data[16,10] = [double]"6.555004832430889"
$data[16,11] = [double]"58.995043491878"
$data[16,12] = [double]"0.03128943259221346"
$data[16,13] = [double]"0.03128943259221347"

# Row 19: MuSCs -> Inflammatory-Mac
$data[17,0] = [double]"0.7761303333333333"
$data[17,1] = [double]"2.328391"
$data[17,2] = [double]"0.03399380269149206"
$data[17,3] = [double]"0.03399380269149207"
$data[17,4] = 2
$data[17,5] = [double]"0.6666666666666666"
$data[17,6] = [double]"0.05870933333333334"
$data[17,7] = [double]"0.176128"
$data[17,8] = [double]"0.006398331365399365"
$data[17,9] = [double]"0.006398331365399365"
$data[17,10] = [double]"0.04556609444977777"
$data[17,11] = [double]"0.410094850048"
$data[17,12] = [double]"0.000217503613990171"
$data[17,13] = [double]"0.0002175036139901711"

# Row 20: MuSCs -> MuSCs
$data[18,0] = [double]"0.7761303333333333"
$data[18,1] = [double]"2.328391"
$data[18,2] = [double]"0.03399380269149206"
$data[18,3] = [double]"0.03399380269149207"
$data[18,4] = 3
$data[18,5] = 1
$data[18,6] = [double]"0.4594193333333333"
$data[18,7] = [double]"1.378258"
$data[18,8] = [double]"0.05006899181852175"
$data[18,9] = [double]"0.05006899181852175"
$data[18,10] = [double]"0.3565692803197777"
$data[18,11] = [double]"3.209123522878"
$data[18,12] = [double]"0.001702035428840759"
$data[18,13] = [double]"0.001702035428840759"

# Row 21: MuSCs -> Resolving-Mac
$data[19,0] = [double]"0.7761303333333333"
$data[19,1] = [double]"2.328391"
$data[19,2] = [double]"0.03399380269149206"
$data[19,3] = [double]"0.03399380269149207"
$data[19,4] = 1
$data[19,5] = [double]"0.3333333333333333"
$data[19,6] = [double]"0.003065"
$data[19,7] = [double]"0.009195"
$data[19,8] = [double]"0.0003340335262130221"
$data[19,9] = [double]"0.0003340335262130222"
$data[19,10] = [double]"0.002378839471666666"
$data[19,11] = [double]"0.021409555245"
$data[19,12] = [double]"1.135506978242882e-05"
$data[19,13] = [double]"1.135506978242882e-05"

# Row 22: Resolving-Mac -> ECs
$data[20,0] = [double]"6.522593333333333"
$data[20,1] = [double]"19.56778"
$data[20,2] = [double]"0.2856836555503455"
$data[20,3] = [double]"0.2856836555503455"
$data[20,4] = 2
$data[20,5] = [double]"0.6666666666666666"
$data[20,6] = [double]"0.2087793333333333"
$data[20,7] = [double]"0.626338"
$data[20,8] = [double]"0.02275344108115409"
$data[20,9] = [double]"0.02275344108115409"
$data[20,10] = [double]"1.361782687737778"
$data[20,11] = [double]"12.25604418964"
$data[20,12] = [double]"0.006500286224413506"
$data[20,13] = [double]"0.006500286224413506"

# Row 23: Resolving-Mac -> FAPs
$data[21,0] = [double]"6.522593333333333"
$data[21,1] = [double]"19.56778"
$data[21,2] = [double]"0.2856836555503455"
$data[21,3] = [double]"0.2856836555503455"
$data[21,4] = 3
$data[21,5] = 1
$data[21,6] = [double]"8.445752666666667"
$data[21,7] = [double]"25.337258"
$data[21,8] = [double]"0.9204452022087118"
$data[21,9] = [double]"0.920445202208712"
$data[21,10] = [double]"55.08821003858223"
$data[21,11] = [double]"495.79389034724"
$data[21,12] = [double]"0.2629561501007618"
$data[21,13] = [double]"0.2629561501007618"

# Row 24: Resolving-Mac -> Inflammatory-Mac
$data[22,0] = [double]"6.522593333333333"
$data[22,1] = [double]"19.56778"
$data[22,2] = [double]"0.2856836555503455"
$data[22,3] = [double]"0.2856836555503455"
$data[22,4] = 2
$data[22,5] = [double]"0.6666666666666666"
$data[22,6] = [double]"0.05870933333333334"
$data[22,7] = [double]"0.176128"
$data[22,8] = [double]"0.006398331365399365"
$data[22,9] = [double]"0.006398331365399365"
$data[22,10] = [double]"0.3829371062044444"
$data[22,11] = [double]"3.44643395584"
$data[22,12] = [double]"0.001827898693889724"
$data[22,13] = [double]"0.001827898693889724"

# Row 25: Resolving-Mac -> MuSCs
$data[23,0] = [double]"6.522593333333333"
$data[23,1] = [double]"19.56778"
$data[23,2] = [double]"0.2856836555503455"
$data[23,3] = [double]"0.2856836555503455"
$data[23,4] = 3
$data[23,5] = 1
$data[23,6] = [double]"0.4594193333333333"
$data[23,7] = [double]"1.378258"
$data[23,8] = [double]"0.05006899181852175"
$data[23,9] = [double]"0.05006899181852175"
$data[23,10] = [double]"2.996605480804444"
$data[23,11] = [double]"26.96944932724"
$data[23,12] = [double]"0.01430389261243564"
$data[23,13] = [double]"0.01430389261243564"

# Row 26: Resolving-Mac -> Resolving-Mac
$data[24,0] = [double]"6.522593333333333"
$data[24,1] = [double]"19.56778"
$data[24,2] = [double]"0.2856836555503455"
$data[24,3] = [double]"0.2856836555503455"
$data[24,4] = 1
$data[24,5] = [double]"0.3333333333333333"
$data[24,6] = [double]"0.003065"
$data[24,7] = [double]"0.009195"
$data[24,8] = [double]"0.0003340335262130221"
$data[24,9] = [double]"0.0003340335262130222"
$data[24,10] = [double]"0.01999174856666667"
$data[24,11] = [double]"0.1799257371"
$data[24,12] = [double]"9.542791884490833e-05"
$data[24,13] = [double]"9.542791884490834e-05"

$ws.Range("G2:T26").Value = $data

Write-Output "Updated G2:T26 with new TPM-derived statistics."
